$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33 was missing the "Project Title" value in column E; shift the
# Month/Year and Amount values one column to the right (E->F, F->G) and
# fill E33 with the project/funding text (duplicate of D33).
$ws.Range("G33").Value = $ws.Range("F33").Value2
$ws.Range("F33").Value = $ws.Range("E33").Value2
$ws.Range("E33").Value = $ws.Range("D33").Value2

# Add the new row 36 entry.
$ws.Range("A36").Value = 2025
$ws.Range("B36").Value = "Dr. Gunjan Mehta"
$ws.Range("E36").Value = "DST-JSPS bilateral grant for organizing an international conference cum workshop on Single-Molecule Biophysics at IIT Hyderabad"
$ws.Range("F36").Value = "2025-2026"
$ws.Range("D36").Value = "DST JSPS"
$ws.Range("H36").Value = "DST"

$ws.Range("G33").Select()
